$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: AD1 = Wins, AE1 = Losses, AF1 = Ties
# Columns: AD = 30, AE = 31, AF = 32
$headers = @("Wins", "Losses", "Ties")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 30 + $i
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Season record for every player row (2 through 41): 83 Wins, 79 Losses, 0 Ties
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 30).Value = 83
    $ws.Cells.Item($r, 31).Value = 79
    $ws.Cells.Item($r, 32).Value = 0
}
